$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values. For Price cells whose new value would otherwise be
# parsed by Excel as a genuine number (losing formatting such as a trailing
# zero, e.g. "16.30" -> 16.3), force the cell to Text format first so the
# original string representation is preserved exactly, matching the
# inline-string cells already used in the workbook.
$ws.Range('D2').Value = '66.985.77'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '3.121.57'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.42'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.53'
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.39'
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.482'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.15'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').Value = '3.638.10'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '66.932.90'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.16'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '3.121.10'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.30'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.44'
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.709'
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E22').Value = '  +4.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.89'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.21'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.30'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.60'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').Value = '0.0₃0951'
$ws.Range('E33').Value = '  -6.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.84'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.975'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '46.91'
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.312'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.60'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').Value = '2.824.50'
$ws.Range('E43').Value = '  +1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '382.71'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('E46').Value = '  -9.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.88'
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.97'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -0.70%  '
